$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Estrela stats update ---
$ws.Range("C12").Value = 23
$ws.Range("F12").Value = 10
$ws.Range("H12").Value = 41
$ws.Range("I12").Value = -15
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = "L D W L L"

# --- Rows 17/18: Tondela and Santa Clara swap places (ranks 16/17) ---
# Row 17 becomes Tondela (previously row 18) with updated stats
$ws.Range("B17").Value = "Tondela"
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 36
$ws.Range("I17").Value = -19
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = 0.78
$ws.Range("L17").Value = "L D D D W"
$ws.Range("M17").Value = 2120
$ws.Range("N17").Value = "Maranhão - 6"
$ws.Range("O17").Value = "Bernardo Fontes"

# Row 18 becomes Santa Clara (previously row 17) with unchanged stats
$ws.Range("B18").Value = "Santa Clara"
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 28
$ws.Range("I18").Value = -11
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = 0.77
$ws.Range("L18").Value = "L L L L L"
$ws.Range("M18").Value = 2920
$ws.Range("N18").Value = "Vinícius - 5"
$ws.Range("O18").Value = "Gabriel Batista"
